$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 209-210 (weekly update adds a new pair of
# Primera/Segunda "Betarraga" observations), shifting all following
# rows down by two.
$ws.Rows("209:210").Insert()

# Row 209 - Calidad "Primera"
$ws.Cells.Item(209, 1).Value = 1
$ws.Cells.Item(209, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(209, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(209, 4).Value = 44642
$ws.Cells.Item(209, 5).Value = 15
$ws.Cells.Item(209, 6).Value = 100114014
$ws.Cells.Item(209, 7).Value = "Betarraga"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 800
$ws.Cells.Item(209, 11).Value = 450
$ws.Cells.Item(209, 12).Value = 500
$ws.Cells.Item(209, 13).Value = 475
$ws.Cells.Item(209, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(209, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(209, 16).Value = 119
$ws.Cells.Item(209, 17).Value = 4
$ws.Cells.Item(209, 18).Value = "Hortaliza"

# Row 210 - Calidad "Segunda"
$ws.Cells.Item(210, 1).Value = 1
$ws.Cells.Item(210, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(210, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(210, 4).Value = 44642
$ws.Cells.Item(210, 5).Value = 15
$ws.Cells.Item(210, 6).Value = 100114014
$ws.Cells.Item(210, 7).Value = "Betarraga"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Segunda"
$ws.Cells.Item(210, 10).Value = 1000
$ws.Cells.Item(210, 11).Value = 450
$ws.Cells.Item(210, 12).Value = 500
$ws.Cells.Item(210, 13).Value = 475
$ws.Cells.Item(210, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(210, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(210, 16).Value = 95
$ws.Cells.Item(210, 17).Value = 5
$ws.Cells.Item(210, 18).Value = "Hortaliza"
